# Update the "Holdings" sheet rows 6-44 (Asset/Industry columns) so they
# carry the real per-row asset/industry labels instead of the placeholder
# blank shared string that was previously repeated on every row.
# (mirrors the already-correct data on the "Constraints" sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Holdings")

$ws.Range("A6").Value = "BND"
$ws.Range("B6").Value = "bonds "
$ws.Range("A7").Value = "GOVT"
$ws.Range("B7").Value = "bonds "
$ws.Range("A8").Value = "SHV"
$ws.Range("B8").Value = "bonds "
$ws.Range("A9").Value = "TIP"
$ws.Range("B9").Value = "bonds "
$ws.Range("A10").Value = "ITE"
$ws.Range("B10").Value = "bonds "
$ws.Range("A11").Value = "MUB"
$ws.Range("B11").Value = "bonds "
$ws.Range("A12").Value = "LQD"
$ws.Range("B12").Value = "bonds "
$ws.Range("A13").Value = "MUB"
$ws.Range("B13").Value = "bonds "
$ws.Range("A14").Value = "LQD"
$ws.Range("B14").Value = "income_generating_bonds "
$ws.Range("A15").Value = "HYG"
$ws.Range("B15").Value = "income_generating_bonds "
$ws.Range("A16").Value = "GOVT"
$ws.Range("B16").Value = "income_generating_bonds "
$ws.Range("A17").Value = "SHV"
$ws.Range("B17").Value = "income_generating_bonds "
$ws.Range("A18").Value = "TIP"
$ws.Range("B18").Value = "income_generating_bonds "
$ws.Range("A19").Value = "ITE"
$ws.Range("B19").Value = "income_generating_bonds "
$ws.Range("A20").Value = "MUB"
$ws.Range("B20").Value = "income_generating_bonds "
$ws.Range("A21").Value = "LQD"
$ws.Range("B21").Value = "income_generating_bonds "
$ws.Range("A22").Value = "MUB"
$ws.Range("B22").Value = "income_generating_bonds "
$ws.Range("A23").Value = "GOVT"
$ws.Range("B23").Value = "income_generating_bonds "
$ws.Range("A24").Value = "GLD"
$ws.Range("B24").Value = "commodities"
$ws.Range("A25").Value = "USO"
$ws.Range("B25").Value = "commodities"
$ws.Range("A26").Value = "DBA"
$ws.Range("B26").Value = "commodities"
$ws.Range("A27").Value = "SLV"
$ws.Range("B27").Value = "commodities"
$ws.Range("A28").Value = "XME"
$ws.Range("B28").Value = "commodities"
$ws.Range("A29").Value = "UNG"
$ws.Range("B29").Value = "commodities"
$ws.Range("A30").Value = "DBB"
$ws.Range("B30").Value = "commodities"
$ws.Range("A31").Value = "GSG"
$ws.Range("B31").Value = "commodities"
$ws.Range("A32").Value = "BNO"
$ws.Range("B32").Value = "commodities"
$ws.Range("A33").Value = "XOP"
$ws.Range("B33").Value = "commodities"
$ws.Range("A34").Value = "BIL"
$ws.Range("B34").Value = "cash_equivalents"
$ws.Range("A35").Value = "SHV"
$ws.Range("B35").Value = "cash_equivalents"
$ws.Range("A36").Value = "TFLO"
$ws.Range("B36").Value = "cash_equivalents"
$ws.Range("A37").Value = "VTI"
$ws.Range("B37").Value = "Equities"
$ws.Range("A38").Value = "IWM"
$ws.Range("B38").Value = "Equities"
$ws.Range("A39").Value = "QQQ"
$ws.Range("B39").Value = "Equities"
$ws.Range("A40").Value = "EEM"
$ws.Range("B40").Value = "Equities"
$ws.Range("A41").Value = "DIA"
$ws.Range("B41").Value = "Equities"
$ws.Range("A42").Value = "IJR"
$ws.Range("B42").Value = "Equities"
$ws.Range("A43").Value = "IVE"
$ws.Range("B43").Value = "Equities"
$ws.Range("A44").Value = "ACWI"
$ws.Range("B44").Value = "Equities"
